# Fruta / hortaliza, semanal
# Insert a new weekly price record row for "Coliflor" (Feria Lagunitas de Puerto Montt)
# at row 345, shifting the existing rows 345:359 down to 346:360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 345; this pushes rows 345-359 down to 346-360
# and Excel copies the formatting (incl. the date number format on column D) from the
# row above, matching the rest of the data table.
$ws.Rows("345:345").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A345").Value = 4
$ws.Range("B345").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C345").Value = "Los Lagos"
$ws.Range("D345").Value = 44753
$ws.Range("E345").Value = 10
$ws.Range("F345").Value = 100112008
$ws.Range("G345").Value = "Coliflor"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 500
$ws.Range("K345").Value = 1800
$ws.Range("L345").Value = 1800
$ws.Range("M345").Value = 1800
$ws.Range("N345").Value = "`$/unidad"
$ws.Range("O345").Value = "Región del Maule"
$ws.Range("P345").Value = 1800
$ws.Range("Q345").Value = 1
$ws.Range("R345").Value = "Hortaliza"
